$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- short-url column (B) changed for every data row ---
$ws.Range("B2:B236").Value = "8EOBan"

# --- refugees (N) / asylum_seekers (O) numeric-as-text updates ---
# Row 228 (item 227): only O changes
$ws.Range("O228").Value = "'30"

# Row 229 (item 228): both N and O change
$ws.Range("N229").Value = "'1129"
$ws.Range("O229").Value = "'394"

# Row 230 (item 229): both N and O change
$ws.Range("N230").Value = "'5695"
$ws.Range("O230").Value = "'8101"

# Row 233 (item 232): both N and O change
$ws.Range("N233").Value = "'13560"
$ws.Range("O233").Value = "'101"

# Row 234 (item 233): both N and O change
$ws.Range("N234").Value = "'17"
$ws.Range("O234").Value = "'6"

# Row 235 (item 234): only O changes
$ws.Range("O235").Value = "'17"

# Row 236 (item 235): only N changes
$ws.Range("N236").Value = "'3553"
